$d = $word.ActiveDocument

$replacements = @(
    @('2025-09-18 Thursday', '2025-09-19 Friday'),
    @('94×63=', '33×63='),
    @('11×45=', '53×28='),
    @('17×71=', '63×54='),
    @('42×73=', '36×74='),
    @('50×69=', '25×79='),
    @('45×51=', '28×19='),
    @('81×90=', '71×81='),
    @('70×15=', '95×55='),
    @('89×33=', '89×58='),
    @('35×89=', '18×28='),
    @('76×69=', '39×40='),
    @('23×41=', '64×75='),
    @('20×11=', '83×72='),
    @('76×51=', '51×47='),
    @('55×74=', '17×93='),
    @('89×40=', '61×77='),
    @('59×86=', '29×97='),
    @('80×19=', '85×48='),
    @('33×95=', '47×96='),
    @('78×54=', '48×92='),
    @('60×28=', '52×84='),
    @('46×33=', '90×97='),
    @('45×29=', '70×39='),
    @('70×24=', '54×23='),
    @('14×58=', '85×71='),
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

$d.Save()
